$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Specify Peter as captain": append ", " + an italic "captain"
#    run to the end of the "Peter Hartman (NetID: pehartma)" item.
# ---------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pehartma*") {
        # $p.Range.End sits just past this paragraph's mark (i.e. at
        # the start of the next paragraph), so back up one character
        # to land right after the closing ")" but still inside this
        # paragraph.
        $paraEnd = $p.Range.End

        $insComma = $d.Range($paraEnd - 1, $paraEnd - 1)
        $insComma.InsertAfter(", ")

        $afterComma = $p.Range.End
        $insCaptain = $d.Range($afterComma - 1, $afterComma - 1)
        $insCaptain.InsertAfter("captain")

        $captainStart = $afterComma - 1
        $captainEnd = $captainStart + 7
        $rCaptain = $d.Range($captainStart, $captainEnd)
        $rCaptain.Font.Italic = $true
        break
    }
}

# ---------------------------------------------------------------
# 2) Remove "(Python)" from the languages answer.
# ---------------------------------------------------------------
$d.Content.Find.Execute("JavaScript, (Python)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "JavaScript", 2)

# ---------------------------------------------------------------
# 3) Add Steven Hernandez's contributions to the workload table
#    (row 3: Task, Estimated Hours, Total).
# ---------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Cell(3, 2).Range.Text = "Parsing module, computational module, validation coordination"
$t.Cell(3, 3).Range.Text = "20"
$t.Cell(3, 4).Range.Text = "20"
